$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" (strikeouts) column values -- replacing the old
# "Strike#" derived numbers with the newly calculated/simulated s_vals.
$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 2
    32 = 0
    33 = 1
    34 = 3
    35 = 1
    36 = 1
    37 = 0
    38 = 0
    40 = 2
    41 = 1
    42 = 0
    43 = 2
    45 = 1
    46 = 1
    47 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
